$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 67.666664
$ws.Range("I11").Value = 67.666664
$ws.Range("K11").Value = 67.666664
$ws.Range("M11").Value = 72.333336
$ws.Range("H17").Value = 2098.889
$ws.Range("J17").Value = 2098.889
$ws.Range("L17").Value = 6296.667
$ws.Range("N17").Value = -6632.667
$ws.Range("H28").Value = 1000
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H32").Value = 3999.9
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 3999.9
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 3999.9
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4651.9
$ws.Range("H48").Value = 2603.45
$ws.Range("I48").Value = 2379.375
$ws.Range("K48").Value = 7138.125
$ws.Range("M48").Value = -6846.125
$ws.Range("H56").Value = 2603.45
$ws.Range("I56").Value = 2379.375
$ws.Range("K56").Value = 7138.125
$ws.Range("M56").Value = -6604.125
$ws.Range("H132").Value = 3236
$ws.Range("I132").Value = 2429.4443
$ws.Range("K132").Value = 7288.3329
$ws.Range("M132").Value = -4758.3329
$ws.Range("H141").Value = 1155.4
$ws.Range("I141").Value = 1155.4
$ws.Range("K141").Value = 3466.2
$ws.Range("M141").Value = 1713.8

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 106662.664
$ws.Range("J125").Value = 106662.664
$ws.Range("L125").Value = 106662.664
$ws.Range("N125").Value = -116502.664
$ws.Range("H132").Value = 33295.184
$ws.Range("I132").Value = 44735.293
$ws.Range("J132").Value = 2788.2222
$ws.Range("K132").Value = 134205.879
$ws.Range("L132").Value = 8364.6666
$ws.Range("M132").Value = -131675.879
$ws.Range("N132").Value = -13424.6666

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8335702.5
$ws.Range("I20").Value = 14287329
$ws.Range("J20").Value = 3425.4
$ws.Range("K20").Value = 14287329
$ws.Range("L20").Value = 3425.4
$ws.Range("M20").Value = -14287082
$ws.Range("N20").Value = -3919.4
$ws.Range("H60").Value = 69894.5
$ws.Range("J60").Value = 69894.5
$ws.Range("L60").Value = 69894.5
$ws.Range("N60").Value = -71092.5
$ws.Range("H86").Value = 2810.8333
$ws.Range("I86").Value = 2856.4
$ws.Range("K86").Value = 2856.4
$ws.Range("M86").Value = -1733.4
$ws.Range("H89").Value = 2810.8333
$ws.Range("I89").Value = 2856.4
$ws.Range("K89").Value = 14282
$ws.Range("M89").Value = -8666
$ws.Range("H99").Value = 49563.617
$ws.Range("I99").Value = 68441.734
$ws.Range("K99").Value = 68441.734
$ws.Range("M99").Value = -66943.734
$ws.Range("H134").Value = 2408.147
$ws.Range("I134").Value = 2393.2424
$ws.Range("J134").Value = 2900
$ws.Range("K134").Value = 7179.7272
$ws.Range("L134").Value = 8700
$ws.Range("M134").Value = -4644.7272
$ws.Range("N134").Value = -13770

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 8499.75
$ws.Range("I42").Value = 8499.75
$ws.Range("K42").Value = 8499.75
$ws.Range("M42").Value = -7906.75
$ws.Range("H58").Value = 23630.436
$ws.Range("J58").Value = 2517.5557
$ws.Range("L58").Value = 2517.5557
$ws.Range("N58").Value = -2923.5557
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 2229.3928
$ws.Range("I132").Value = 1711.8667
$ws.Range("J132").Value = 2826.5386
$ws.Range("K132").Value = 5135.6001
$ws.Range("L132").Value = 8479.6158
$ws.Range("M132").Value = -2605.6001
$ws.Range("N132").Value = -13539.6158
$ws.Range("H136").Value = 23630.436
$ws.Range("J136").Value = 2517.5557
$ws.Range("L136").Value = 7552.6671
$ws.Range("N136").Value = -12652.6671
$ws.Range("H137").Value = 72142.14
$ws.Range("J137").Value = 99998.75
$ws.Range("L137").Value = 99998.75
$ws.Range("N137").Value = -110198.75

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 287.5
$ws.Range("I23").Value = 283.33334
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 850.0000200000001
$ws.Range("L23").Value = 900
$ws.Range("M23").Value = -615.0000200000001
$ws.Range("N23").Value = -1370
$ws.Range("H34").Value = 968.8461
$ws.Range("I34").Value = 966.25
$ws.Range("K34").Value = 2898.75
$ws.Range("M34").Value = -2814.75
$ws.Range("H41").Value = 453.7037
$ws.Range("I41").Value = 450
$ws.Range("J41").Value = 550
$ws.Range("K41").Value = 1350
$ws.Range("L41").Value = 1650
$ws.Range("M41").Value = -1012
$ws.Range("N41").Value = -2326

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H53").Value = 49999
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H70").Value = 8332.333000000001
$ws.Range("I70").Value = 9999
$ws.Range("K70").Value = 9999
$ws.Range("M70").Value = -9729
$ws.Range("H73").Value = 8332.333000000001
$ws.Range("I73").Value = 9999
$ws.Range("K73").Value = 9999
$ws.Range("M73").Value = -9063
$ws.Range("H80").Value = 3612.7144
$ws.Range("I80").Value = 3497.25
$ws.Range("K80").Value = 3497.25
$ws.Range("M80").Value = -2499.25
$ws.Range("H83").Value = 3612.7144
$ws.Range("I83").Value = 3497.25
$ws.Range("K83").Value = 17486.25
$ws.Range("M83").Value = -12494.25
$ws.Range("H97").Value = 1358
$ws.Range("I97").Value = 310.66666
$ws.Range("K97").Value = 310.66666
$ws.Range("M97").Value = 185.33334
$ws.Range("H132").Value = 33187.242
$ws.Range("J132").Value = 2955.5
$ws.Range("L132").Value = 8866.5
$ws.Range("N132").Value = -13926.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3520.1177
$ws.Range("I16").Value = 3677.6875
$ws.Range("J16").Value = 999
$ws.Range("K16").Value = 3677.6875
$ws.Range("L16").Value = 999
$ws.Range("M16").Value = -3507.6875
$ws.Range("N16").Value = -1339
$ws.Range("H43").Value = 27824.5
$ws.Range("J43").Value = 33413.43
$ws.Range("L43").Value = 33413.43
$ws.Range("N43").Value = -33799.43

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 23760.5
$ws.Range("I52").Value = 25347.334
$ws.Range("K52").Value = 25347.334
$ws.Range("M52").Value = -25121.334
$ws.Range("H54").Value = 3452.3076
$ws.Range("I54").Value = 3452.3076
$ws.Range("K54").Value = 3452.3076
$ws.Range("M54").Value = -2932.3076
$ws.Range("H58").Value = 21514.166
$ws.Range("I58").Value = 21514.166
$ws.Range("K58").Value = 21514.166
$ws.Range("M58").Value = -21206.166
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
